{"js": "// Update the two-digit multiplication problem/answer strings throughout the\n// document's tables. Each entry is the exact original cell text (a unique\n// \"A\u00d7B=C\" string) paired with its replacement, matching the canonical OOXML\n// diff exactly.\nconst replacements = [\n  [\"72\u00d729=2088\", \"72\u00d753=3816\"],\n  [\"93\u00d758=5394\", \"55\u00d731=1705\"],\n  [\"79\u00d713=1027\", \"86\u00d789=7654\"],\n  [\"76\u00d772=5472\", \"58\u00d715=870\"],\n  [\"80\u00d750=4000\", \"96\u00d723=2208\"],\n  [\"43\u00d769=2967\", \"30\u00d777=2310\"],\n  [\"58\u00d772=4176\", \"19\u00d713=247\"],\n  [\"47\u00d757=2679\", \"48\u00d764=3072\"],\n  [\"36\u00d783=2988\", \"70\u00d749=3430\"],\n  [\"71\u00d743=3053\", \"73\u00d790=6570\"],\n  [\"45\u00d711=495\", \"67\u00d748=3216\"],\n  [\"69\u00d772=4968\", \"55\u00d738=2090\"],\n  [\"96\u00d719=1824\", \"57\u00d761=3477\"],\n  [\"91\u00d717=1547\", \"52\u00d711=572\"],\n  [\"88\u00d726=2288\", \"37\u00d774=2738\"],\n  [\"34\u00d758=1972\", \"24\u00d777=1848\"],\n  [\"97\u00d713=1261\", \"78\u00d769=5382\"],\n  [\"99\u00d789=8811\", \"35\u00d723=805\"],\n  [\"98\u00d793=9114\", \"75\u00d749=3675\"],\n  [\"99\u00d777=7623\", \"83\u00d791=7553\"],\n  [\"48\u00d759=2832\", \"18\u00d787=1566\"],\n  [\"19\u00d784=1596\", \"16\u00d789=1424\"],\n  [\"42\u00d797=4074\", \"92\u00d741=3772\"],\n  [\"55\u00d791=5005\", \"49\u00d779=3871\"],\n  [\"59\u00d729=1711\", \"72\u00d760=4320\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problem/answer strings throughout the\n# document's tables. Each entry is the exact original cell text (a unique\n# \"A\u00d7B=C\" string) paired with its replacement, matching the canonical OOXML\n# diff exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"72\u00d729=2088\", \"72\u00d753=3816\"),\n  @(\"93\u00d758=5394\", \"55\u00d731=1705\"),\n  @(\"79\u00d713=1027\", \"86\u00d789=7654\"),\n  @(\"76\u00d772=5472\", \"58\u00d715=870\"),\n  @(\"80\u00d750=4000\", \"96\u00d723=2208\"),\n  @(\"43\u00d769=2967\", \"30\u00d777=2310\"),\n  @(\"58\u00d772=4176\", \"19\u00d713=247\"),\n  @(\"47\u00d757=2679\", \"48\u00d764=3072\"),\n  @(\"36\u00d783=2988\", \"70\u00d749=3430\"),\n  @(\"71\u00d743=3053\", \"73\u00d790=6570\"),\n  @(\"45\u00d711=495\", \"67\u00d748=3216\"),\n  @(\"69\u00d772=4968\", \"55\u00d738=2090\"),\n  @(\"96\u00d719=1824\", \"57\u00d761=3477\"),\n  @(\"91\u00d717=1547\", \"52\u00d711=572\"),\n  @(\"88\u00d726=2288\", \"37\u00d774=2738\"),\n  @(\"34\u00d758=1972\", \"24\u00d777=1848\"),\n  @(\"97\u00d713=1261\", \"78\u00d769=5382\"),\n  @(\"99\u00d789=8811\", \"35\u00d723=805\"),\n  @(\"98\u00d793=9114\", \"75\u00d749=3675\"),\n  @(\"99\u00d777=7623\", \"83\u00d791=7553\"),\n  @(\"48\u00d759=2832\", \"18\u00d787=1566\"),\n  @(\"19\u00d784=1596\", \"16\u00d789=1424\"),\n  @(\"42\u00d797=4074\", \"92\u00d741=3772\"),\n  @(\"55\u00d791=5005\", \"49\u00d779=3871\"),\n  @(\"59\u00d729=1711\", \"72\u00d760=4320\"),\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute([ref]$old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
